$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.164.71"
$ws.Range("D3").Value = "2.321.23"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.33%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").Value = "2.682.89"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "2.314.61"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.797"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "43.101.96"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("D21").Value = "0.0₃0909"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  +3.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.78"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.58%  "
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0698"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D42").Value = "2.000.42"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0289"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.79%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "2.548.04"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("E51").Value = "  +1.65%  "

Write-Host "done"
